# Added Larissa the Priest
# Row 13 on the "Characters" sheet already has Larissa/F/Priest set up;
# this fills in her weapon and stat line (the row previously had zeroes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Characters")

# Weapon: Staff -> Spirits (new shared string)
$ws.Range("D13").Value = "Spirits"

# Stat block (Health, MP, Attack, Defense, Resistance, Skill, Speed)
$ws.Range("G13").Value = 50
$ws.Range("H13").Value = 70
$ws.Range("I13").Value = 60
$ws.Range("J13").Value = 40
$ws.Range("K13").Value = 70
$ws.Range("L13").Value = 40
$ws.Range("M13").Value = 55

# Matches the post-edit selection left on the sheet
$null = $ws.Range("L14").Select()
